# "4 Jan Presenti Sheet" - fill in attendance for 4-Jan column (G) on the
# Jan-2024 sheet, mirroring F5's comment pattern, and widen the data
# validation / note so G participates like C:F already do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Attendance values for 4-Jan (column G) ---------------------------
$ws.Range("G2").Value = "Absent"
$ws.Range("G3").Value = "Present"
$ws.Range("G4").Value = "Present"
$ws.Range("G5").Value = "Present"

# --- Comment on G2 explaining the Absent mark --------------------------
$comment = $ws.Range("G2").AddComment()
[void]$comment.Text("A:" + [char]10 + "university Exam")

# --- Extend the Present/Absent/Reason validation list to include G -----
$ws.Range("C2:G5").Validation.Delete()
$ws.Range("C2:G5").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')
$ws.Range("C2:G5").Validation.IgnoreBlank = $true
$ws.Range("C2:G5").Validation.InCellDropdown = $true
$ws.Range("C2:G5").Validation.ShowInput = $true
$ws.Range("C2:G5").Validation.ShowError = $true

# --- Move the active selection the way the author left it --------------
[void]$ws.Range("H16").Select()
